$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add border around AA11 (top+left+right thin) - existing formula cell
# ---------------------------------------------------------------------------
$r = $ws.Range("AA11")
$r.Borders.Item(7).LineStyle = 1   ' left
$r.Borders.Item(7).Weight = 2
$r.Borders.Item(10).LineStyle = 1  ' right
$r.Borders.Item(10).Weight = 2
$r.Borders.Item(8).LineStyle = 1   ' top
$r.Borders.Item(8).Weight = 2

# ---------------------------------------------------------------------------
# 2. New labels in column Z (rows 12-15) - plain text labels
# ---------------------------------------------------------------------------
$ws.Range("Z12").Value = "y10"
$ws.Range("Z13").Value = "L21"
$ws.Range("Z14").Value = "L20"
$ws.Range("Z15").Value = "S3"

# ---------------------------------------------------------------------------
# 3. New labels in column Z (rows 16-19) - label with superscript "2"
# ---------------------------------------------------------------------------
$ws.Range("Z16").Value = "y102"
$ws.Range("Z16").Characters(4, 1).Font.Superscript = $true

$ws.Range("Z17").Value = "L212"
$ws.Range("Z17").Characters(4, 1).Font.Superscript = $true

$ws.Range("Z18").Value = "L202"
$ws.Range("Z18").Characters(4, 1).Font.Superscript = $true

$ws.Range("Z19").Value = "S32"
$ws.Range("Z19").Characters(3, 1).Font.Superscript = $true

# ---------------------------------------------------------------------------
# 4. New formulas in column AA (rows 12-19), with borders + number formats
#    matching the source cell they echo.
# ---------------------------------------------------------------------------

# AA12 = Y10  (left/right thin border only, general format)
$ws.Range("AA12").Formula = "=Y10"
$c = $ws.Range("AA12")
$c.Borders.Item(7).LineStyle = 1
$c.Borders.Item(7).Weight = 2
$c.Borders.Item(10).LineStyle = 1
$c.Borders.Item(10).Weight = 2

# AA13 = L21  (left/right thin border, number format 0.00000)
$ws.Range("AA13").Formula = "=L21"
$c = $ws.Range("AA13")
$c.Borders.Item(7).LineStyle = 1
$c.Borders.Item(7).Weight = 2
$c.Borders.Item(10).LineStyle = 1
$c.Borders.Item(10).Weight = 2
$c.NumberFormat = "0.00000"

# AA14 = L20  (left/right thin border, number format 0.000)
$ws.Range("AA14").Formula = "=L20"
$c = $ws.Range("AA14")
$c.Borders.Item(7).LineStyle = 1
$c.Borders.Item(7).Weight = 2
$c.Borders.Item(10).LineStyle = 1
$c.Borders.Item(10).Weight = 2
$c.NumberFormat = "0.000"

# AA15 = S3  (left/right thin border, number format 0.0000)
$ws.Range("AA15").Formula = "=S3"
$c = $ws.Range("AA15")
$c.Borders.Item(7).LineStyle = 1
$c.Borders.Item(7).Weight = 2
$c.Borders.Item(10).LineStyle = 1
$c.Borders.Item(10).Weight = 2
$c.NumberFormat = "0.0000"

# AA16 = AA12^2  (left/right thin border only, general format)
$ws.Range("AA16").Formula = "=AA12^2"
$c = $ws.Range("AA16")
$c.Borders.Item(7).LineStyle = 1
$c.Borders.Item(7).Weight = 2
$c.Borders.Item(10).LineStyle = 1
$c.Borders.Item(10).Weight = 2

# AA17 = AA13^2  (left/right thin border only, general format)
$ws.Range("AA17").Formula = "=AA13^2"
$c = $ws.Range("AA17")
$c.Borders.Item(7).LineStyle = 1
$c.Borders.Item(7).Weight = 2
$c.Borders.Item(10).LineStyle = 1
$c.Borders.Item(10).Weight = 2

# AA18 = AA14^2  (left/right thin border only, general format)
$ws.Range("AA18").Formula = "=AA14^2"
$c = $ws.Range("AA18")
$c.Borders.Item(7).LineStyle = 1
$c.Borders.Item(7).Weight = 2
$c.Borders.Item(10).LineStyle = 1
$c.Borders.Item(10).Weight = 2

# AA19 = AA15^2  (left/right/bottom thin border, general format)
$ws.Range("AA19").Formula = "=AA15^2"
$c = $ws.Range("AA19")
$c.Borders.Item(7).LineStyle = 1
$c.Borders.Item(7).Weight = 2
$c.Borders.Item(10).LineStyle = 1
$c.Borders.Item(10).Weight = 2
$c.Borders.Item(9).LineStyle = 1
$c.Borders.Item(9).Weight = 2

# ---------------------------------------------------------------------------
# 5. Row height adjustments for rows 16-19 (taller rows to fit superscript)
# ---------------------------------------------------------------------------
$ws.Rows(16).RowHeight = 16.5
$ws.Rows(17).RowHeight = 16.5
$ws.Rows(18).RowHeight = 16.5
$ws.Rows(19).RowHeight = 17

# ---------------------------------------------------------------------------
# 6. Update the sheet view: scroll so column I is the left-most visible
#    column, and move the active selection to Z20.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 9
$ws.Range("Z20").Select()
